$d = $word.ActiveDocument

# Replace the title "Clase 1 – ReactJS" with "Módulo 2 – Estilos y Event Handlers"
$d.Content.Find.Execute("Clase 1 – ReactJS", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Módulo 2 – Estilos y Event Handlers", 2)
